$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vcan"
$ws.Range("C2").Value = "Cd44"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.094072
$ws.Range("H2").Value = 12.282216
$ws.Range("I2").Value = 0.0147134235951136
$ws.Range("J2").Value = 0.0147134235951136
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 239.0839323333333
$ws.Range("N2").Value = 717.251797
$ws.Range("O2").Value = 0.4086975387666237
$ws.Range("P2").Value = 0.4086975387666237
$ws.Range("Q2").Value = 978.8268330157946
$ws.Range("R2").Value = 8809.441497142152
$ws.Range("S2").Value = 0.006013340010153695
$ws.Range("T2").Value = 0.006013340010153695

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vcan"
$ws.Range("C3").Value = "Cd44"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.094072
$ws.Range("H3").Value = 12.282216
$ws.Range("I3").Value = 0.0147134235951136
$ws.Range("J3").Value = 0.0147134235951136
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 117.0512696666667
$ws.Range("N3").Value = 351.153809
$ws.Range("O3").Value = 0.2000910950200451
$ws.Range("P3").Value = 0.2000910950200451
$ws.Range("Q3").Value = 479.2163257067493
$ws.Range("R3").Value = 4312.946931360744
$ws.Range("S3").Value = 0.002944025038640048
$ws.Range("T3").Value = 0.002944025038640048

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vcan"
$ws.Range("C4").Value = "Cd44"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 4.094072
$ws.Range("H4").Value = 12.282216
$ws.Range("I4").Value = 0.0147134235951136
$ws.Range("J4").Value = 0.0147134235951136
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 171.15883
$ws.Range("N4").Value = 513.47649
$ws.Range("O4").Value = 0.2925842480357353
$ws.Range("P4").Value = 0.2925842480357353
$ws.Range("Q4").Value = 700.73657345576
$ws.Range("R4").Value = 6306.629161101841
$ws.Range("S4").Value = 0.004304915978607556
$ws.Range("T4").Value = 0.004304915978607557

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Vcan"
$ws.Range("C5").Value = "Cd44"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 4.094072
$ws.Range("H5").Value = 12.282216
$ws.Range("I5").Value = 0.0147134235951136
$ws.Range("J5").Value = 0.0147134235951136
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 57.695868
$ws.Range("N5").Value = 173.087604
$ws.Range("O5").Value = 0.09862711817759588
$ws.Range("P5").Value = 0.09862711817759588
$ws.Range("Q5").Value = 236.211037694496
$ws.Range("R5").Value = 2125.899339250464
$ws.Range("S5").Value = 0.001451142567712296
$ws.Range("T5").Value = 0.001451142567712297

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vcan"
$ws.Range("C6").Value = "Cd44"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 181.0215506666667
$ws.Range("H6").Value = 543.064652
$ws.Range("I6").Value = 0.6505617768331834
$ws.Range("J6").Value = 0.6505617768331835
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 239.0839323333333
$ws.Range("N6").Value = 717.251797
$ws.Range("O6").Value = 0.4086975387666237
$ws.Range("P6").Value = 0.4086975387666237
$ws.Range("Q6").Value = 43279.34417046441
$ws.Range("R6").Value = 389514.0975341797
$ws.Range("S6").Value = 0.2658829970073636
$ws.Range("T6").Value = 0.2658829970073636

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vcan"
$ws.Range("C7").Value = "Cd44"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 181.0215506666667
$ws.Range("H7").Value = 543.064652
$ws.Range("I7").Value = 0.6505617768331834
$ws.Range("J7").Value = 0.6505617768331835
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 117.0512696666667
$ws.Range("N7").Value = 351.153809
$ws.Range("O7").Value = 0.2000910950200451
$ws.Range("P7").Value = 0.2000910950200451
$ws.Range("Q7").Value = 21188.80234256216
$ws.Range("R7").Value = 190699.2210830595
$ws.Range("S7").Value = 0.1301716183047379
$ws.Range("T7").Value = 0.1301716183047379

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Vcan"
$ws.Range("C8").Value = "Cd44"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 181.0215506666667
$ws.Range("H8").Value = 543.064652
$ws.Range("I8").Value = 0.6505617768331834
$ws.Range("J8").Value = 0.6505617768331835
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 171.15883
$ws.Range("N8").Value = 513.47649
$ws.Range("O8").Value = 0.2925842480357353
$ws.Range("P8").Value = 0.2925842480357353
$ws.Range("Q8").Value = 30983.43681689239
$ws.Range("R8").Value = 278850.9313520315
$ws.Range("S8").Value = 0.1903441282755288
$ws.Range("T8").Value = 0.1903441282755288

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Vcan"
$ws.Range("C9").Value = "Cd44"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 181.0215506666667
$ws.Range("H9").Value = 543.064652
$ws.Range("I9").Value = 0.6505617768331834
$ws.Range("J9").Value = 0.6505617768331835
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 57.695868
$ws.Range("N9").Value = 173.087604
$ws.Range("O9").Value = 0.09862711817759588
$ws.Range("P9").Value = 0.09862711817759588
$ws.Range("Q9").Value = 10444.19549241931
$ws.Range("R9").Value = 93997.75943177381
$ws.Range("S9").Value = 0.06416303324555314
$ws.Range("T9").Value = 0.06416303324555314

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Vcan"
$ws.Range("C10").Value = "Cd44"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 63.05609033333334
$ws.Range("H10").Value = 189.168271
$ws.Range("I10").Value = 0.2266132513854375
$ws.Range("J10").Value = 0.2266132513854376
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 239.0839323333333
$ws.Range("N10").Value = 717.251797
$ws.Range("O10").Value = 0.4086975387666237
$ws.Range("P10").Value = 0.4086975387666237
$ws.Range("Q10").Value = 15075.69803445922
$ws.Range("R10").Value = 135681.282310133
$ws.Range("S10").Value = 0.09261627809313049
$ws.Range("T10").Value = 0.09261627809313051

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Vcan"
$ws.Range("C11").Value = "Cd44"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 63.05609033333334
$ws.Range("H11").Value = 189.168271
$ws.Range("I11").Value = 0.2266132513854375
$ws.Range("J11").Value = 0.2266132513854376
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 117.0512696666667
$ws.Range("N11").Value = 351.153809
$ws.Range("O11").Value = 0.2000910950200451
$ws.Range("P11").Value = 0.2000910950200451
$ws.Range("Q11").Value = 7380.795433732694
$ws.Range("R11").Value = 66427.15890359425
$ws.Range("S11").Value = 0.04534329361576494
$ws.Range("T11").Value = 0.04534329361576495

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Vcan"
$ws.Range("C12").Value = "Cd44"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 63.05609033333334
$ws.Range("H12").Value = 189.168271
$ws.Range("I12").Value = 0.2266132513854375
$ws.Range("J12").Value = 0.2266132513854376
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 171.15883
$ws.Range("N12").Value = 513.47649
$ws.Range("O12").Value = 0.2925842480357353
$ws.Range("P12").Value = 0.2925842480357353
$ws.Range("Q12").Value = 10792.60664582764
$ws.Range("R12").Value = 97133.4598124488
$ws.Range("S12").Value = 0.06630346775154129
$ws.Range("T12").Value = 0.06630346775154131

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Vcan"
$ws.Range("C13").Value = "Cd44"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 63.05609033333334
$ws.Range("H13").Value = 189.168271
$ws.Range("I13").Value = 0.2266132513854375
$ws.Range("J13").Value = 0.2266132513854376
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 57.695868
$ws.Range("N13").Value = 173.087604
$ws.Range("O13").Value = 0.09862711817759588
$ws.Range("P13").Value = 0.09862711817759588
$ws.Range("Q13").Value = 3638.075864468076
$ws.Range("R13").Value = 32742.68278021269
$ws.Range("S13").Value = 0.02235021192500079
$ws.Range("T13").Value = 0.0223502119250008

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Vcan"
$ws.Range("C14").Value = "Cd44"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 30.082493
$ws.Range("H14").Value = 90.247479
$ws.Range("I14").Value = 0.1081115481862653
$ws.Range("J14").Value = 0.1081115481862653
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 239.0839323333333
$ws.Range("N14").Value = 717.251797
$ws.Range("O14").Value = 0.4086975387666237
$ws.Range("P14").Value = 0.4086975387666237
$ws.Range("Q14").Value = 7192.240720829974
$ws.Range("R14").Value = 64730.16648746976
$ws.Range("S14").Value = 0.04418492365597587
$ws.Range("T14").Value = 0.04418492365597588

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Vcan"
$ws.Range("C15").Value = "Cd44"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 30.082493
$ws.Range("H15").Value = 90.247479
$ws.Range("I15").Value = 0.1081115481862653
$ws.Range("J15").Value = 0.1081115481862653
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 117.0512696666667
$ws.Range("N15").Value = 351.153809
$ws.Range("O15").Value = 0.2000910950200451
$ws.Range("P15").Value = 0.2000910950200451
$ws.Range("Q15").Value = 3521.194000388612
$ws.Range("R15").Value = 31690.74600349751
$ws.Range("S15").Value = 0.0216321580609022
$ws.Range("T15").Value = 0.0216321580609022

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Vcan"
$ws.Range("C16").Value = "Cd44"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 30.082493
$ws.Range("H16").Value = 90.247479
$ws.Range("I16").Value = 0.1081115481862653
$ws.Range("J16").Value = 0.1081115481862653
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 171.15883
$ws.Range("N16").Value = 513.47649
$ws.Range("O16").Value = 0.2925842480357353
$ws.Range("P16").Value = 0.2925842480357353
$ws.Range("Q16").Value = 5148.88430536319
$ws.Range("R16").Value = 46339.95874826871
$ws.Range("S16").Value = 0.0316317360300576
$ws.Range("T16").Value = 0.0316317360300576

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Vcan"
$ws.Range("C17").Value = "Cd44"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 30.082493
$ws.Range("H17").Value = 90.247479
$ws.Range("I17").Value = 0.1081115481862653
$ws.Range("J17").Value = 0.1081115481862653
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 57.695868
$ws.Range("N17").Value = 173.087604
$ws.Range("O17").Value = 0.09862711817759588
$ws.Range("P17").Value = 0.09862711817759588
$ws.Range("Q17").Value = 1735.635545238924
$ws.Range("R17").Value = 15620.71990715032
$ws.Range("S17").Value = 0.01066273043932964
$ws.Range("T17").Value = 0.01066273043932964
